$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data that was in row 2 and row 3 got swapped (e.g. "Tjäder" / "Spillkråka"
# records exchanged rows), while the "Publik kommentar" note ("2 tuppar") stayed
# attached to the Tjäder record and therefore moved from row 2 to row 3.
#
# Swap every cell that actually differs between row 2 and row 3 (columns that
# are identical in both rows are left untouched so no unrelated formatting or
# values are disturbed).

$cols = @("A","B","D","E","F","G","H","Q","R","Z","AB")
foreach ($col in $cols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $v3
    $ws.Range($addr3).Value2 = $v2
}

# "Publik kommentar" ("2 tuppar") moves from AC2 to AC3 (AC3 was empty before).
$ws.Range("AC3").Value2 = $ws.Range("AC2").Value2
$ws.Range("AC2").ClearContents()
